# The sheet holding the imported "user_5" dataset was renamed from
# "Sheet2" to "Sheet1". Excel automatically keeps the dependent "user_5"
# defined name's reference in sync with the rename
# (Sheet2!$A$1:$G$104 -> Sheet1!$A$1:$G$104).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "Sheet1"

# The workbook also picked up the hidden "LOCAL_MYSQL_DATE_FORMAT" helper
# name that the MySQL for Excel add-in stores in workbooks it has touched.
$mysqlDateFormat = "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)"
$mysqlName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", $mysqlDateFormat)
$mysqlName.Visible = $false
